# Add ANR for electricity production to new technology dataset
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("NewTechFramework")

# New rows of data (PlantType, DataSource, ATBTechnologyType, FuelType,
# ThermalOrRenewableOrStorage, Capacity (MW), ..., NSPSCompliant,
# NOxEmRate, SO2EmRate, CO2EmRate, Lifetime(years), ..., SO2 Scrubber,
# CoalType, Efficiency, ..., Minimum Energy Capacity, Maximum Charge Rate,
# ECAPEX, PlantCategory)
$rows = @(
    @{ Row = 18; PlantType = "iPWR";     Capacity = 77;   Lifetime = 60 },
    @{ Row = 19; PlantType = "HTGR";     Capacity = 164;  Lifetime = 60 },
    @{ Row = 20; PlantType = "PBRHTGR";  Capacity = 80;   Lifetime = 60 },
    @{ Row = 21; PlantType = "iMSR";     Capacity = 141;  Lifetime = 60 },
    @{ Row = 22; PlantType = "Micro";    Capacity = 6.7;  Lifetime = 20 }
)

# Populate column by column (A for all rows, then B for all rows, etc.)
# so that new shared-string entries are introduced in the same order
# as the authored workbook.
foreach ($r in $rows) {
    $ws.Cells.Item($r.Row, 1).Value = $r.PlantType    # A - PlantType
}
foreach ($r in $rows) {
    $ws.Cells.Item($r.Row, 2).Value = "ANRElec"       # B - DataSource
}
foreach ($r in $rows) {
    $ws.Cells.Item($r.Row, 3).Value = "NA"            # C - ATBTechnologyType
}
foreach ($r in $rows) {
    $ws.Cells.Item($r.Row, 4).Value = "Nuclear Fuel"  # D - FuelType
}
foreach ($r in $rows) {
    $ws.Cells.Item($r.Row, 5).Value = "thermal"       # E - ThermalOrRenewableOrStorage
}
foreach ($r in $rows) {
    $ws.Cells.Item($r.Row, 6).Value = $r.Capacity     # F - Capacity (MW)
}
foreach ($r in $rows) {
    $ws.Cells.Item($r.Row, 11).Value = "Yes"          # K - NSPSCompliant
}
foreach ($r in $rows) {
    $ws.Cells.Item($r.Row, 12).Value = 0              # L - NOxEmRate(lb/MMBtu)
}
foreach ($r in $rows) {
    $ws.Cells.Item($r.Row, 13).Value = 0              # M - SO2EmRate(lb/MMBtu)
}
foreach ($r in $rows) {
    $ws.Cells.Item($r.Row, 14).Value = 0              # N - CO2EmRate(lb/MMBtu)
}
foreach ($r in $rows) {
    $ws.Cells.Item($r.Row, 15).Value = $r.Lifetime    # O - Lifetime(years)
}
foreach ($r in $rows) {
    $ws.Cells.Item($r.Row, 17).Value = "NA"           # Q - SO2 Scrubber
}
foreach ($r in $rows) {
    $ws.Cells.Item($r.Row, 18).Value = "NA"           # R - CoalType
}
foreach ($r in $rows) {
    $ws.Cells.Item($r.Row, 19).Value = "NA"           # S - Efficiency
}
foreach ($r in $rows) {
    $ws.Cells.Item($r.Row, 21).Value = "NA"           # U - Minimum Energy Capacity (MWh)
}
foreach ($r in $rows) {
    $ws.Cells.Item($r.Row, 22).Value = "NA"           # V - Maximum Charge Rate (MW)
}
foreach ($r in $rows) {
    $ws.Cells.Item($r.Row, 23).Value = "NA"           # W - ECAPEX(2012$/MWH)
}
foreach ($r in $rows) {
    $ws.Cells.Item($r.Row, 24).Value = "Nuclear"      # X - PlantCategory
}

$ws.Range("X18:X22").Select()
